$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.018.79"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.562.09"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  +0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("E6").Value = "  +0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.08"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0597"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "1.783.81"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D13").Value = "1.562.79"
$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").Value = "27.039.01"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.84"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.85%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.05"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("E28").Value = "  +1.54%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0474"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.37%  "

$ws.Range("E31").Value = "  +3.51%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("E33").Value = "  +3.62%  "

$ws.Range("D34").Value = "1.422.22"
$ws.Range("E34").Value = "  +0.35%  "

$ws.Range("E35").Value = "  +1.97%  "

$ws.Range("E36").Value = "  +10.28%  "

$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.533"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("E40").Value = "  +0.73%  "

$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("E42").Value = "  +0.39%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.18%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("D47").Value = "1.698.47"
$ws.Range("E47").Value = "  +0.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("E49").Value = "  +3.24%  "

$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.77%  "
